$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "592.00", "134.51") are not auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.393.89'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '3.509.93'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '592.00'
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("D6").Value = '134.51'
$ws.Range("E6").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  +5.97%  '

$ws.Range("E10").Value = '  +0.61%  '

$ws.Range("E11").Value = '  +3.30%  '

$ws.Range("D12").Value = '4.111.04'

$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("D15").Value = '3.513.57'
$ws.Range("E15").Value = '  +0.60%  '

$ws.Range("D16").Value = '25.75'
$ws.Range("E16").Value = '  +2.34%  '

$ws.Range("D17").Value = '64.390.67'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("D19").Value = '13.64'
$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("D20").Value = '5.75'
$ws.Range("E20").Value = '  +1.82%  '

$ws.Range("D21").Value = '394.39'
$ws.Range("E21").Value = '  +2.19%  '

$ws.Range("E22").Value = '  +1.91%  '

$ws.Range("D23").Value = '3.650.62'
$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("D24").Value = '74.65'

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("E26").Value = '  +0.45%  '

$ws.Range("E27").Value = '  +3.34%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").Value = '8.25'
$ws.Range("E31").Value = '  +0.49%  '

$ws.Range("D32").Value = '1.48'
$ws.Range("E32").Value = '  -3.58%  '

$ws.Range("E33").Value = '  +6.85%  '

$ws.Range("D34").Value = '3.539.19'
$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").Value = '23.35'
$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("D37").Value = '5.36'
$ws.Range("E37").Value = '  +1.17%  '

$ws.Range("D38").Value = '6.95'
$ws.Range("E38").Value = '  +1.75%  '

$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("D40").Value = '166.92'
$ws.Range("E40").Value = '  +2.71%  '

$ws.Range("D41").Value = '0.0787'
$ws.Range("E41").Value = '  +0.69%  '

$ws.Range("E42").Value = '  +0.65%  '

$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("E44").Value = '  +0.98%  '

$ws.Range("D45").Value = '25.10'
$ws.Range("E45").Value = '  -2.92%  '

$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("E47").Value = '  -3.07%  '

$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").Value = '0.911'
$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.377.60'
$ws.Range("E50").Value = '  -4.05%  '

$ws.Range("E51").Value = '  +0.23%  '

# Restore default cell style (removes the temporary Text number format
# marker) while keeping the values stored as literal text.
$dRange.Style = "Normal"
